# Applies the "several months of updates" backup edit:
#  1. Bumps the cached datetimeFigureOut placeholder text (2021-04-14 -> 2021-12-14)
#     on the slide master and every slide layout.
#  2. Nudges the "data_intergroup table" rectangle (shape id 11) on slide 1.
#  3. Re-routes the two connectors feeding it (shape ids 97 and 100) on slide 1.
#  4. Shortens two description strings on slide 1 (shape ids 21 and 130).

$p = $ppt.ActivePresentation

function Get-ShapeById($shapes, $id) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.Id -eq $id) {
            return $shp
        }
    }
    return $null
}

function Set-DatePlaceholderText($shapes, $newText) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.Name -like "Date Placeholder*") {
            $shp.TextFrame.TextRange.Text = $newText
        }
    }
}

function Replace-Substring($textRange, $oldSubstring, $newSubstring) {
    $full = $textRange.Text
    $idx = $full.IndexOf($oldSubstring)
    if ($idx -ge 0) {
        $chars = $textRange.Characters($idx + 1, $oldSubstring.Length)
        $chars.Text = $newSubstring
    }
}

# --- 1. Update the cached date field everywhere it is rendered -------------

$newDate = "2021-12-14"

$master = $p.SlideMaster
Set-DatePlaceholderText $master.Shapes $newDate

$layouts = $master.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    $layout = $layouts.Item($li)
    Set-DatePlaceholderText $layout.Shapes $newDate
}

# --- 2/3. Re-position the data_intergroup rectangle + its connectors -------

$slide1 = $p.Slides.Item(1)

$dataIntergroup = Get-ShapeById $slide1.Shapes 11
if ($dataIntergroup -ne $null) {
    $dataIntergroup.Left = 666.07016048031493937742
    $dataIntergroup.Top = 32.77118210236220363640
}

$connLine = Get-ShapeById $slide1.Shapes 97
if ($connLine -ne $null) {
    $connLine.Left = 640.21385826771654592449
    $connLine.Top = 112.87945181889764967309
    $connLine.Width = 0.00007874015748031496
    $connLine.Height = 182.72654043307085203196
}

$connArrow = Get-ShapeById $slide1.Shapes 100
if ($connArrow -ne $null) {
    $connArrow.Left = 640.21385826771654592449
    $connArrow.Top = 112.87945181889764967309
}

# --- 4. Trim the two description strings ------------------------------------

$modelGroupsShape = Get-ShapeById $slide1.Shapes 21
if ($modelGroupsShape -ne $null) {
    Replace-Substring $modelGroupsShape.TextFrame.TextRange "functional groups for the ECOSPACE model groups" "functional groups / model groups"
}

$ssSpeciesShape = Get-ShapeById $slide1.Shapes 130
if ($ssSpeciesShape -ne $null) {
    Replace-Substring $ssSpeciesShape.TextFrame.TextRange "species list for Salish Sea (should include verified TSN codes)" "species list (should include verified TSN codes)"
}

Write-Host "edit.ps1 applied"
